$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 586, shifting existing rows 586:683 down to 587:684
$ws.Rows.Item(586).Insert()

# Populate the newly inserted row 586 with the new weekly data record
$ws.Cells.Item(586, 1).Value = 9
$ws.Cells.Item(586, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(586, 3).Value = "Metropolitana"
$ws.Cells.Item(586, 4).Value = 45218
$ws.Cells.Item(586, 5).Value = 13
$ws.Cells.Item(586, 6).Value = 100112012
$ws.Cells.Item(586, 7).Value = "Espinaca"
$ws.Cells.Item(586, 8).Value = "Sin especificar"
$ws.Cells.Item(586, 9).Value = "Primera"
$ws.Cells.Item(586, 10).Value = 52
$ws.Cells.Item(586, 11).Value = 10000
$ws.Cells.Item(586, 12).Value = 12000
$ws.Cells.Item(586, 13).Value = 11000
$ws.Cells.Item(586, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(586, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(586, 16).Value = 1100
$ws.Cells.Item(586, 17).Value = 10
$ws.Cells.Item(586, 18).Value = "Hortaliza"
